$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.342.19'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '1.836.88'
$ws.Range('E3').Value = '  +3.32%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.557'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.03'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +4.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0724'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0931'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('E12').Value = '  +3.39%  '
$ws.Range('D13').Value = '1.837.19'
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.647'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('E15').Value = '  -3.01%  '
$ws.Range('D16').Value = '34.376.28'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.35'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('E20').Value = '  +8.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.31%  '
$ws.Range('E23').Value = '  +2.22%  '
$ws.Range('E24').Value = '  +1.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('E27').Value = '  +4.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.116'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0538'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.09%  '
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('E32').Value = '  +2.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.59'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('D35').Value = '1.451.10'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.650'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.39%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.59%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0193'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.971'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '82.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('E41').Value = '  -2.93%  '
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('E43').Value = '  +4.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.14%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.996.72'
$ws.Range('E45').Value = '  +3.18%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0501'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +8.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.94'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.95%  '
$ws.Range('E51').Value = '  +6.33%  '
